# Update the "Estado de Cuenta" worker data table (rows 16-19) so that
# each worker's periods of arrears ("Periodo Mora") are grouped together.
# DAMARIS ESTHER BENEDETTY BONFANTE now occupies rows 16-17 (periods 1701,1612)
# and EDITH MARIA ROMERO MARTINEZ now occupies rows 18-19 (periods 1701,1612).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - keep worker, change period
$ws.Range("C16").Value = "1047365859"
$ws.Range("D16").Value = "DAMARIS ESTHER BENEDETTY BONFANTE"
$ws.Range("E16").Value = "1701"

# Row 17 - switch worker, keep period
$ws.Range("C17").Value = "1047365859"
$ws.Range("D17").Value = "DAMARIS ESTHER BENEDETTY BONFANTE"
$ws.Range("E17").Value = "1612"

# Row 18 - unchanged, set explicitly for consistency
$ws.Range("C18").Value = "1143379924"
$ws.Range("D18").Value = "EDITH MARIA ROMERO MARTINEZ"
$ws.Range("E18").Value = "1701"

# Row 19 - switch worker, change period
$ws.Range("C19").Value = "1143379924"
$ws.Range("D19").Value = "EDITH MARIA ROMERO MARTINEZ"
$ws.Range("E19").Value = "1612"
